$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 2 — opportunity 1329970 -> 1329997 (Recruitment Consultant, Prague)
# ---------------------------------------------------------------------------
$ws.Cells.Item(2,2).Value = "https://aiesec.org/opportunity/global-talent/1329997"
$ws.Cells.Item(2,3).Value = "Recruitment Consultant"
$ws.Cells.Item(2,4).Value = "Prague, Czechia"
$ws.Cells.Item(2,6).Value = "8 applicants"
$ws.Cells.Item(2,8).Value = "Non Stop Consulting"

# ---------------------------------------------------------------------------
# Row 3 — opportunity 1329856 -> 1327511 (IT Research Interns (Duplicated))
# ---------------------------------------------------------------------------
$ws.Cells.Item(3,2).Value = "https://aiesec.org/opportunity/global-talent/1327511"
$ws.Cells.Item(3,3).Value = "IT Research Interns (Duplicated)"
$ws.Cells.Item(3,4).Value = "Aronj, Uttar Pradesh, India"
# E3 loses its special "Yes" highlight style and becomes a plain "No"
$ws.Cells.Item(3,5).Style = "Normal"
$ws.Cells.Item(3,5).Value = "No"
$ws.Cells.Item(3,6).Value = "7 applicants"
$ws.Cells.Item(3,7).Value = "3 - 6 Months"
$ws.Cells.Item(3,8).Value = "FS University"

# ---------------------------------------------------------------------------
# Row 4 — opportunity 1329697 -> 1325379 (Software Development Intern)
# ---------------------------------------------------------------------------
$ws.Cells.Item(4,2).Value = "https://aiesec.org/opportunity/global-talent/1325379"
$ws.Cells.Item(4,3).Value = "Software Development Intern"
$ws.Cells.Item(4,4).Value = "Athens, Greece"
$ws.Cells.Item(4,6).Value = "129 applicants"
$ws.Cells.Item(4,8).Value = "Eutopians"

# ---------------------------------------------------------------------------
# Row 5 — opportunity 1325604 -> 1316660 (IT Research Interns)
# ---------------------------------------------------------------------------
$ws.Cells.Item(5,2).Value = "https://aiesec.org/opportunity/global-talent/1316660"
$ws.Cells.Item(5,3).Value = "IT Research Interns"
$ws.Cells.Item(5,4).Value = "Aronj, Uttar Pradesh, India"
$ws.Cells.Item(5,6).Value = "6 applicants"
$ws.Cells.Item(5,8).Value = "FS University"

# ---------------------------------------------------------------------------
# Column A holds opportunity IDs that look numeric ("1329997" etc.) but must
# stay stored as text, exactly like the rest of the sheet. Force text entry
# mode via a temporary "@" number format, then drop back to the Normal style
# so the cells end up with no explicit style (matching the original file).
# ---------------------------------------------------------------------------
$ws.Range("A2:A5").NumberFormat = "@"
$ws.Cells.Item(2,1).Value = "1329997"
$ws.Cells.Item(3,1).Value = "1327511"
$ws.Cells.Item(4,1).Value = "1325379"
$ws.Cells.Item(5,1).Value = "1316660"
$ws.Range("A2:A5").Style = "Normal"

# ---------------------------------------------------------------------------
# Row 6 (opportunity 1316099) is dropped entirely.
# ---------------------------------------------------------------------------
$ws.Rows.Item(6).Delete()

# ---------------------------------------------------------------------------
# Column width tweaks for C, D and H (ColumnWidth uses character units, which
# store as width+0.8333 in the sheet XML, so subtract that padding here).
# ---------------------------------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 34.16666666666667
$ws.Columns.Item(4).ColumnWidth = 29.166666666666668
$ws.Columns.Item(8).ColumnWidth = 21.166666666666668

Write-Output "edit complete"
